# Apply text replacements to the multiplication problems table.
$d = $word.ActiveDocument

$replacements = @(
    @{old = "338×6="; new = "436×4="},
    @{old = "973×3="; new = "385×6="},
    @{old = "494×5="; new = "765×2="},
    @{old = "920×6="; new = "283×2="},
    @{old = "283×3="; new = "319×4="},
    @{old = "571×5="; new = "659×5="},
    @{old = "339×3="; new = "492×6="},
    @{old = "673×9="; new = "238×7="},
    @{old = "135×3="; new = "590×6="},
    @{old = "869×5="; new = "131×5="},
    @{old = "893×6="; new = "367×2="},
    @{old = "383×9="; new = "819×5="},
    @{old = "447×8="; new = "719×4="},
    @{old = "805×5="; new = "182×6="},
    @{old = "197×6="; new = "564×9="},
    @{old = "490×6="; new = "462×9="},
    @{old = "890×5="; new = "381×4="},
    @{old = "293×3="; new = "189×2="},
    @{old = "836×5="; new = "563×5="},
    @{old = "775×3="; new = "345×2="},
    @{old = "714×9="; new = "892×8="},
    @{old = "657×9="; new = "136×2="},
    @{old = "904×5="; new = "516×5="},
    @{old = "870×4="; new = "943×2="},
    @{old = "140×6="; new = "801×8="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $true, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
